# Apply the package-republish update to StructureDefinition-ConsentPeriod.xlsx
#
# Sheet 1 ("Metadata"): bump Version/Status/Date/FHIR Version, and insert a new
#   "Jurisdiction" property row right after "Contact".
# Sheet 2 ("Elements"): three text corrections (Constraint(s) on the root
#   Extension row, Type(s) on Extension.id, and the R4B->R4 link fix on
#   Extension.value[x]).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet ---------------------------------------------------

$meta.Range("B3").Value  = "2.0.2"                                 # Version
$meta.Range("B6").Value  = "active"                                # Status
$meta.Range("B8").Value  = "2025-02-05T10:42:38+00:00"             # Date

# Insert a new "Jurisdiction" row right after the "Contact" row (row 10),
# copying its formatting so the new row keeps the same bordered style.
$meta.Rows.Item(11).Insert()
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# FHIR Version now lives one row further down because of the inserted row.
$meta.Range("B15").Value = "4.0.1"

# --- Elements sheet -----------------------------------------------------

# Row 2 = root "Extension" element: Constraint(s) / AJ2
$elements.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 3 = "Extension.id": Type(s) / K3
$elements.Range("K3").Value = "string" + [char]10

# Row 6 = "Extension.value[x]": Definition / M6
$elements.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
